# Apply cryptos list price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.351.31"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").Value = "3.969.31"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "613.46"
$ws.Range("E5").Value = "  +1.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.67"
$ws.Range("E6").Value = "  +3.98%  "
$ws.Range("D7").Value = "3.967.64"
$ws.Range("E7").Value = "  +2.94%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("E9").Value = "  +1.12%  "
$ws.Range("E10").Value = "  +2.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("E11").Value = "  +2.65%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("E13").Value = "  +5.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.37"
$ws.Range("E14").Value = "  +3.39%  "
$ws.Range("D15").Value = "4.623.65"
$ws.Range("E15").Value = "  +2.64%  "
$ws.Range("D16").Value = "3.954.79"
$ws.Range("E16").Value = "  +2.21%  "
$ws.Range("D17").Value = "70.258.97"
$ws.Range("E17").Value = "  +1.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.66"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.90"
$ws.Range("E19").Value = "  +4.14%  "
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.14"
$ws.Range("E21").Value = "  -4.46%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "503.60"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.744"
$ws.Range("E23").Value = "  +3.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000169"
$ws.Range("E24").Value = "  +6.63%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.80"
$ws.Range("E25").Value = "  +1.43%  "
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.50"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("E28").Value = "  +2.98%  "
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").Value = "4.119.61"
$ws.Range("E31").Value = "  +2.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.44"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.94"
$ws.Range("E33").Value = "  -0.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "32.57"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "3.934.39"
$ws.Range("E35").Value = "  +3.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.109"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.20"
$ws.Range("E37").Value = "  +4.97%  "
$ws.Range("E38").Value = "  +0.76%  "
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.31"
$ws.Range("E40").Value = "  +9.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").Value = "  -0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.328"
$ws.Range("E42").Value = "  +2.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.08"
$ws.Range("E43").Value = "  +4.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "441.73"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "48.35"
$ws.Range("E45").Value = "  -0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.67"
$ws.Range("E46").Value = "  +2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000278"
$ws.Range("E48").Value = "  +22.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0369"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.41"
$ws.Range("E50").Value = "  +2.41%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "143.40"
$ws.Range("E51").Value = "  +0.21%  "
